# Fix last name issue in fellows: remove the stray row for "teredolz@cotaipec.org.mx"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the cell holding the incorrect/duplicate email address and delete its entire row
$found = $ws.UsedRange.Find("teredolz@cotaipec.org.mx")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}

# Update the selected cell to match the post-edit state
$ws.Range("C42").Select()
